# Automatische test-sync: 2025-06-26 23:24:50
# Appends a new test-mail row (row 39) to the "Logs" sheet and bumps the
# matching category counter on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$row = 39

$logs.Cells.Item($row, 1).Value = "Zijn jullie telefonisch bereikbaar?"
$logs.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($row, 3).Value = "Testmail #7: Zijn jullie telefonisch bereikbaar?"
$logs.Cells.Item($row, 4).Value = "Openingstijden / Locatie"
$logs.Cells.Item($row, 5).Value = "Beste klant,`nBedankt voor uw e-mail. Op dit moment zijn wij helaas niet telefonisch bereikbaar. Als u vragen heeft of hulp nodig heeft, kunt u ons altijd per e-mail bereiken op dit adres. We streven ernaar om binnen 24 uur te reageren op alle vragen en verzoeken die we ontvangen.`nMet vriendelijke groet,`n[Naam van het bedrijf] assistent"
$logs.Cells.Item($row, 6).Value = "2025-06-26 23:24:20"
$logs.Cells.Item($row, 7).Value = "Ja"
$logs.Cells.Item($row, 8).Value = "Nee"
$logs.Cells.Item($row, 9).Value = "Ja"

# Dashboard: "Openingstijden / Locatie" count goes from 10 to 11
$dashboard.Cells.Item(3, 2).Value = 11
